# Insert a new data row before row 33 (Fruta / hortaliza, semanal), pushing
# the existing rows 33-47 down to 34-48 and extending the table by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 45093
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100107
$ws.Cells.Item(33, 8).Value = "Otros"
$ws.Cells.Item(33, 9).Value = 100107011
$ws.Cells.Item(33, 10).Value = "Tuna"
$ws.Cells.Item(33, 11).Value = "Sin especificar"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 40
$ws.Cells.Item(33, 14).Value = 22000
$ws.Cells.Item(33, 15).Value = 22000
$ws.Cells.Item(33, 16).Value = 22000
$ws.Cells.Item(33, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(33, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(33, 19).Value = 1222
$ws.Cells.Item(33, 20).Value = 18
